$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dollar counts for Los Yahoo (row 3) and Out of PO (row 4)
$ws.Range("B3").Value = 401
$ws.Range("B4").Value = 403

# Add a new ledger entry in row 13 (matching the plain-text style used by row 11)
$ws.Range("A13").Value = "21.02.2025 - Out of PO GTJ karşılığında (Wemby'i unutma bedeli) Los Yahoo'ya 1 Dolar vermiştir. (403-401)"
$ws.Range("A11").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Update the ledger note text in A12 (Portis / dize yatırma entry reworded)
$ws.Range("A12").Value = "19.02.2025 - Otistics Portis karşılığında (dize yatırma eylemi Portisin 25 maç ceza almasıyla tersine döndü) Los Yahoo'ya 1 Dolar vermiştir. (394-400)"

# Leave the cursor where the author ended up after editing
$ws.Range("B19").Select() | Out-Null
